# Daily attendance processing - 2026-01-12 22:59:42
# Normalizes the "Recorded By" (column G) audit-trail strings for rows whose
# recorder list currently places "System" ahead of the other recorder(s).
# The most-recent recorder is rotated to the front of the comma-separated
# list (equivalently, "System" is moved to the end of its entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after mappings observed for the "Recorded By" column.
$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
